$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 89156.336
$ws.Range("J47").Value = 89156.336
$ws.Range("L47").Value = 89156.336
$ws.Range("N47").Value = -91100.336
$ws.Range("H86").Value = 69447930
$ws.Range("I86").Value = 41670210
$ws.Range("J86").Value = 83336780
$ws.Range("K86").Value = 41670210
$ws.Range("L86").Value = 83336780
$ws.Range("M86").Value = -41669087
$ws.Range("N86").Value = -83339026
$ws.Range("H89").Value = 69447930
$ws.Range("I89").Value = 41670210
$ws.Range("J89").Value = 83336780
$ws.Range("K89").Value = 208351050
$ws.Range("L89").Value = 416683900
$ws.Range("M89").Value = -208345434
$ws.Range("N89").Value = -416695132
$ws.Range("H96").Value = 3386
$ws.Range("I96").Value = 2228
$ws.Range("J96").Value = 10334
$ws.Range("K96").Value = 6684
$ws.Range("L96").Value = 31002
$ws.Range("M96").Value = -5311
$ws.Range("N96").Value = -33748
$ws.Range("H99").Value = 502.53845
$ws.Range("J99").Value = 984.2
$ws.Range("L99").Value = 2952.6
$ws.Range("N99").Value = -5948.6
$ws.Range("H107").Value = 144.33333
$ws.Range("I107").Value = 181.42857
$ws.Range("K107").Value = 181.42857
$ws.Range("M107").Value = 1738.57143
$ws.Range("H109").Value = 46481.3
$ws.Range("J109").Value = 46481.3
$ws.Range("L109").Value = 46481.3
$ws.Range("N109").Value = -49255.3
$ws.Range("H114").Value = 99989.336
$ws.Range("J114").Value = 99989.336
$ws.Range("L114").Value = 99989.336
$ws.Range("N114").Value = -108667.336
$ws.Range("H132").Value = 1048.431
$ws.Range("I132").Value = 1126.5098
$ws.Range("J132").Value = 479.57144
$ws.Range("K132").Value = 3379.5294
$ws.Range("L132").Value = 1438.71432
$ws.Range("M132").Value = -849.5294000000004
$ws.Range("N132").Value = -6498.71432
$ws.Range("H138").Value = 1443.825
$ws.Range("I138").Value = 959.7778
$ws.Range("J138").Value = 2449.1538
$ws.Range("K138").Value = 2879.3334
$ws.Range("L138").Value = 7347.4614
$ws.Range("M138").Value = 2260.6666
$ws.Range("N138").Value = -17627.4614
$ws.Range("H141").Value = 8773.666999999999
$ws.Range("I141").Value = 8160.6665
$ws.Range("J141").Value = 9999.666999999999
$ws.Range("K141").Value = 24481.9995
$ws.Range("L141").Value = 29999.001
$ws.Range("M141").Value = -19301.9995
$ws.Range("N141").Value = -40359.001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 1824.5
$ws.Range("I36").Value = 1824.5
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1824.5
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1478.5
$ws.Range("N36").ClearContents()
$ws.Range("H61").Value = 1513.3846
$ws.Range("I61").Value = 1278.4762
$ws.Range("K61").Value = 1278.4762
$ws.Range("M61").Value = -1066.4762
$ws.Range("H96").Value = 35000
$ws.Range("J96").Value = 35000
$ws.Range("L96").Value = 35000
$ws.Range("N96").Value = -40492
$ws.Range("H97").Value = 744
$ws.Range("I97").Value = 744
$ws.Range("K97").Value = 744
$ws.Range("M97").Value = -248
$ws.Range("H132").Value = 1681.3636
$ws.Range("I132").Value = 1344.8182
$ws.Range("K132").Value = 4034.4546
$ws.Range("M132").Value = -1504.4546
$ws.Range("H136").Value = 1513.3846
$ws.Range("I136").Value = 1278.4762
$ws.Range("K136").Value = 3835.4286
$ws.Range("M136").Value = -1285.4286

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 12144.714
$ws.Range("I86").Value = 8335.333000000001
$ws.Range("J86").Value = 15001.75
$ws.Range("K86").Value = 8335.333000000001
$ws.Range("L86").Value = 15001.75
$ws.Range("M86").Value = -7212.333000000001
$ws.Range("N86").Value = -17247.75
$ws.Range("H89").Value = 12144.714
$ws.Range("I89").Value = 8335.333000000001
$ws.Range("J89").Value = 15001.75
$ws.Range("K89").Value = 41676.665
$ws.Range("L89").Value = 75008.75
$ws.Range("M89").Value = -36060.665
$ws.Range("N89").Value = -86240.75
$ws.Range("H94").Value = 3240
$ws.Range("I94").Value = 3240
$ws.Range("K94").Value = 3240
$ws.Range("M94").Value = -2789
$ws.Range("H134").Value = 5227.091
$ws.Range("J134").Value = 7214.143
$ws.Range("L134").Value = 21642.429
$ws.Range("N134").Value = -26712.429

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2088.1875
$ws.Range("I31").Value = 1257.9286
$ws.Range("K31").Value = 1257.9286
$ws.Range("M31").Value = -962.9286
$ws.Range("H34").Value = 2088.1875
$ws.Range("I34").Value = 1257.9286
$ws.Range("K34").Value = 1257.9286
$ws.Range("M34").Value = -1055.9286
$ws.Range("H58").Value = 1626.6666
$ws.Range("I58").Value = 1424.7222
$ws.Range("K58").Value = 1424.7222
$ws.Range("M58").Value = -1221.7222
$ws.Range("H99").Value = 7409792
$ws.Range("I99").Value = 11113120
$ws.Range("K99").Value = 11113120
$ws.Range("M99").Value = -11111622
$ws.Range("H105").Value = 5000
$ws.Range("I105").Value = 2500
$ws.Range("J105").Value = 10000
$ws.Range("K105").Value = 2500
$ws.Range("L105").Value = 10000
$ws.Range("M105").Value = -753
$ws.Range("N105").Value = -13494
$ws.Range("H122").Value = 2246.9092
$ws.Range("J122").Value = 2834.25
$ws.Range("L122").Value = 8502.75
$ws.Range("N122").Value = -13402.75
$ws.Range("H126").Value = 7409792
$ws.Range("I126").Value = 11113120
$ws.Range("K126").Value = 33339360
$ws.Range("M126").Value = -33336890
$ws.Range("H132").Value = 1792.4783
$ws.Range("I132").Value = 1641.8462
$ws.Range("J132").Value = 1988.3
$ws.Range("K132").Value = 4925.5386
$ws.Range("L132").Value = 5964.9
$ws.Range("M132").Value = -2395.5386
$ws.Range("N132").Value = -11024.9
$ws.Range("H136").Value = 1626.6666
$ws.Range("I136").Value = 1424.7222
$ws.Range("K136").Value = 4274.1666
$ws.Range("M136").Value = -1724.1666

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 9349.666999999999
$ws.Range("I18").Value = 10819.6
$ws.Range("K18").Value = 32458.8
$ws.Range("M18").Value = -32289.8
$ws.Range("H23").Value = 83467.836
$ws.Range("I23").Value = 15
$ws.Range("J23").Value = 91054.45
$ws.Range("K23").Value = 45
$ws.Range("L23").Value = 273163.35
$ws.Range("N23").Value = -273633.35
$ws.Range("M23").Value = 190
$ws.Range("H39").Value = 6533.0625
$ws.Range("I39").Value = 300
$ws.Range("K39").Value = 900
$ws.Range("M39").Value = -606
$ws.Range("H50").Value = 244.5
$ws.Range("I50").Value = 244.5
$ws.Range("K50").Value = 733.5
$ws.Range("M50").Value = -252.5
$ws.Range("H53").Value = 244.5
$ws.Range("I53").Value = 244.5
$ws.Range("K53").Value = 733.5
$ws.Range("M53").Value = -252.5
$ws.Range("H110").Value = 4750
$ws.Range("I110").Value = 4750
$ws.Range("K110").Value = 14250
$ws.Range("M110").Value = -10160
$ws.Range("H111").Value = 375
$ws.Range("I111").Value = 375
$ws.Range("K111").Value = 1125
$ws.Range("M111").Value = 1942

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H57").Value = 16999.666
$ws.Range("J57").Value = 40000
$ws.Range("L57").Value = 40000
$ws.Range("N57").Value = -41640
$ws.Range("H132").Value = 2851.8462
$ws.Range("J132").Value = 5058.9287
$ws.Range("L132").Value = 15176.7861
$ws.Range("N132").Value = -20236.7861

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1762.75
$ws.Range("I22").Value = 1462.5
$ws.Range("K22").Value = 1462.5
$ws.Range("M22").Value = -1167.5
$ws.Range("H27").Value = 1762.75
$ws.Range("I27").Value = 1462.5
$ws.Range("K27").Value = 1462.5
$ws.Range("M27").Value = -1355.5
$ws.Range("H55").Value = 2997.4443
$ws.Range("I55").Value = 371.16666
$ws.Range("K55").Value = 371.16666
$ws.Range("M55").Value = -198.16666
$ws.Range("H132").Value = 2493.5217
$ws.Range("I132").Value = 2214.1667
$ws.Range("K132").Value = 6642.500100000001
$ws.Range("M132").Value = -4112.500100000001
$ws.Range("H136").Value = 2656.3572
$ws.Range("I136").Value = 3872.5833
$ws.Range("K136").Value = 11617.7499
$ws.Range("M136").Value = -9067.749899999999
